$wb = $excel.ActiveWorkbook

$wsBeth  = $wb.Worksheets.Item("Beth S.")
$wsChris = $wb.Worksheets.Item("Chris K.")
$wsMary  = $wb.Worksheets.Item("Mary M.")

# --- Mary M. (sheet3): add her shopping-list items ---
$wsMary.Range("A2").Value = "Pizza"
$wsMary.Range("B2").Value = 1
$wsMary.Range("A3").Value = "Fruits "
$wsMary.Range("B3").Value = 1
$wsMary.Range("A4").Value = "Vegetables"
$wsMary.Range("B4").Value = 1
$wsMary.Range("A5").Value = "Basic Dog Food"
$wsMary.Range("B5").Value = 1
[void]$wsMary.Range("B6").Select()

# --- Chris K. (sheet2): add his shopping-list items ---
$wsChris.Range("A2").Value = "Cat Litter"
$wsChris.Range("B2").Value = 1
$wsChris.Range("A3").Value = "Premium Cat Food"
$wsChris.Range("B3").Value = 1
$wsChris.Range("A4").Value = "Brush"
$wsChris.Range("B4").Value = 1
$wsChris.Range("A5").Value = "Oatmeal Soap"
$wsChris.Range("B5").Value = 1
[void]$wsChris.Range("B6").Select()

# --- Beth S. (sheet1): add her shopping-list items ---
$wsBeth.Range("A2").Value = "Body Butter"
$wsBeth.Range("B2").Value = 1
$wsBeth.Range("A3").Value = "Catnip"
$wsBeth.Range("B3").Value = 1
$wsBeth.Range("A4").Value = "Fruits"
$wsBeth.Range("B4").Value = 1
$wsBeth.Range("A5").Value = "Vegetables"
$wsBeth.Range("B5").Value = 1
[void]$wsBeth.Range("B6").Select()
